$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (dependent SUM formulas recalc automatically)
$ws.Range("N9").Value = 1
$ws.Range("I14").Value = 1

# Apply underline font formatting to E14 and G15
$ws.Range("E14").Font.Underline = $true
$ws.Range("G15").Font.Underline = $true

# Update selection / active cell
$null = $ws.Range("N8").Select()
